# Apply the edits described by the commit:
#  1. On slide 4, merge the two runs in the second body paragraph into a
#     single run (the sentence about "additional bathrooms" loses its
#     separate-run formatting break).
#  2. Append a new slide 5 ("What the data doesn't show") using the same
#     "Title and Content" layout as the rest of the deck, with a title
#     that has an italic "doesn't" and a two-paragraph body.

$p = $ppt.ActivePresentation

# --- 1. Fix slide 4 ------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$body4 = $slide4.Shapes.Item(2).TextFrame.TextRange
$para2 = $body4.Paragraphs(2, 1)
# A no-op assignment of text identical to what's already there leaves the
# existing run split untouched, so first change it to something else and
# then set the real (merged) text; that forces the paragraph back down to
# a single run.
$para2.Text = "temp"
$para2.Text = "Following square footage, additional bedrooms are more strongly correlated with assessed value than additional bathrooms"

# --- 2. Add slide 5 -------------------------------------------------------
$titleAndContent = $p.SlideMaster.CustomLayouts.Item(2)
$slide5 = $p.Slides.AddSlide($p.Slides.Count + 1, $titleAndContent)

$apostrophe = [char]0x2019

# Title: "What the data doesn't show" with "doesn't" italicised.
$title5 = $slide5.Shapes.Item(1).TextFrame.TextRange
$title5.Text = "What the data doesn" + $apostrophe + "t show"
$title5.Characters(15, 7).Font.Italic = $true

# Body: the placeholder's "no-bullet" margins (marL=0 / indent=0) can only
# be pushed through TextFrame2.Ruler onto paragraph 1, and InsertBefore
# inherits the paragraph formatting of the paragraph it's attached to. So
# build the three paragraphs back-to-front: start with paragraph 3's
# (unformatted) text, prepend paragraph 2 and fix its formatting while it
# is paragraph 1, then prepend paragraph 1 and let it inherit that fix.

$shape5 = $slide5.Shapes.Item(2)
$body5 = $shape5.TextFrame.TextRange

$body5.Text = "Forecasting from this data would be greatly improved by merging in Census and Labor Department data and by using additional data from the database to establish local metrics."

$body5.InsertBefore("Los Angeles County includes 87 different cities plus dozens of other legal entities. These include Beverly Hills, Malibu, Compton, Watts, and Beautiful Downtown Burbank. These all have vastly different factors driving pricing that are more significant that anything captured in the data set. The data should be localized.`r")

$ruler5 = $shape5.TextFrame2.Ruler
$ruler5.Levels.Item(1).LeftMargin = 0
$ruler5.Levels.Item(1).FirstMargin = 0
$body5.Paragraphs(1, 1).ParagraphFormat.Bullet.Type = 0

$body5.InsertBefore("Regardless of stated preferences, this dataset only covers Los Angeles County.`r")
$ruler5b = $shape5.TextFrame2.Ruler
$ruler5b.Levels.Item(1).LeftMargin = 0
$ruler5b.Levels.Item(1).FirstMargin = 0
$body5.Paragraphs(1, 1).ParagraphFormat.Bullet.Type = 0

Write-Host "Slide count now:" $p.Slides.Count
